# Adding the subjects table: the old sequential "row number" helper column
# (column C, values 1..49) is removed from Sheet1, and the last user
# selection is now the whole of column C (as if the user clicked the
# column header to select it before deleting its contents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded helper column C (was just =ROW() style 1..49 values)
$ws.Range("C1:C49").ClearContents()

# Reflect the user's final selection: the entire column C
$ws.Columns("C").Select()
